# Auto-generated edit script applying scheduled market-data refresh
# to Atomos_Profits leve-profit columns (H-N) across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 6 (Leve Item ID 4564)
$ws.Range("H6").Value = 1681.5385
$ws.Range("I6").Value = 160
$ws.Range("K6").Value = 480
$ws.Range("M6").Value = -368

# ALC row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 703.625
$ws.Range("I33").Value = 701.2308
$ws.Range("J33").Value = 706.4545000000001
$ws.Range("K33").Value = 701.2308
$ws.Range("L33").Value = 706.4545000000001
$ws.Range("M33").Value = -472.2308
$ws.Range("N33").Value = -1164.4545

# ALC row 42 (Leve Item ID 4600)
$ws.Range("H42").Value = 1338
$ws.Range("I42").Value = 230
$ws.Range("J42").Value = 3000
$ws.Range("K42").Value = 690
$ws.Range("L42").Value = 9000
$ws.Range("M42").Value = -460
$ws.Range("N42").Value = -9460

# ALC row 112 (Leve Item ID 27960)
$ws.Range("H112").Value = 1530.1818
$ws.Range("J112").Value = 1145.2222
$ws.Range("L112").Value = 3435.6666
$ws.Range("N112").Value = -5651.6666

# ALC row 135 (Leve Item ID 44047)
$ws.Range("H135").Value = 904.7778
$ws.Range("I135").Value = 710.3913
$ws.Range("J135").Value = 2022.5
$ws.Range("K135").Value = 6393.5217
$ws.Range("L135").Value = 18202.5
$ws.Range("M135").Value = -3858.5217
$ws.Range("N135").Value = -23272.5

# ALC row 141 (Leve Item ID 44161)
$ws.Range("H141").Value = 714830.7
$ws.Range("I141").Value = 2279.182
$ws.Range("J141").Value = 1834554.4
$ws.Range("K141").Value = 6837.545999999999
$ws.Range("L141").Value = 5503663.199999999
$ws.Range("M141").Value = -1657.545999999999
$ws.Range("N141").Value = -5514023.199999999

$ws = $wb.Worksheets.Item("ARM")
# ARM row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 1515.8387
$ws.Range("I45").Value = 1116.7084
$ws.Range("J45").Value = 2884.2856
$ws.Range("K45").Value = 1116.7084
$ws.Range("L45").Value = 2884.2856
$ws.Range("M45").Value = -739.7084
$ws.Range("N45").Value = -3638.2856

# ARM row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 1709.1333
$ws.Range("I61").Value = 818.2308
$ws.Range("J61").Value = 7500
$ws.Range("K61").Value = 818.2308
$ws.Range("L61").Value = 7500
$ws.Range("M61").Value = -606.2308
$ws.Range("N61").Value = -7924

# ARM row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 1709.1333
$ws.Range("I136").Value = 818.2308
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 2454.6924
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = 95.30760000000009
$ws.Range("N136").Value = -27600

$ws = $wb.Worksheets.Item("BSM")
# BSM row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 3999.6667
$ws.Range("I107").Value = 899.4
$ws.Range("K107").Value = 899.4
$ws.Range("M107").Value = 1020.6

$ws = $wb.Worksheets.Item("CRP")
# CRP row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 1349.6666
$ws.Range("I107").Value = 542.55554
$ws.Range("K107").Value = 542.55554
$ws.Range("M107").Value = 1377.44446

# CRP row 141 (Leve Item ID 43345)
$ws.Range("H141").Value = 20202.166
$ws.Range("J141").Value = 20202.166
$ws.Range("L141").Value = 20202.166
$ws.Range("N141").Value = -30562.166

$ws = $wb.Worksheets.Item("CUL")
# CUL row 10 (Leve Item ID 4689)
$ws.Range("H10").Value = 956.2
$ws.Range("I10").Value = 156
$ws.Range("J10").Value = 1756.4
$ws.Range("K10").Value = 468
$ws.Range("L10").Value = 5269.200000000001
$ws.Range("M10").Value = -329
$ws.Range("N10").Value = -5547.200000000001

# CUL row 35 (Leve Item ID 4718)
$ws.Range("H35").Value = 2911.8572
$ws.Range("J35").Value = 3576.6
$ws.Range("L35").Value = 10729.8
$ws.Range("N35").Value = -11305.8

# CUL row 47 (Leve Item ID 4663)
$ws.Range("H47").Value = 2498.1428
$ws.Range("I47").Value = 291.5
$ws.Range("K47").Value = 874.5
$ws.Range("M47").Value = -443.5

# CUL row 123 (Leve Item ID 36037)
$ws.Range("H123").Value = 3049.0908
$ws.Range("I123").Value = 513.3333
$ws.Range("J123").Value = 4000
$ws.Range("K123").Value = 1539.9999
$ws.Range("L123").Value = 12000
$ws.Range("M123").Value = 910.0001
$ws.Range("N123").Value = -16900

# CUL row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 1325.8971
$ws.Range("J131").Value = 1112.2203
$ws.Range("L131").Value = 3336.6609
$ws.Range("N131").Value = -13416.6609

# CUL row 134 (Leve Item ID 44074)
$ws.Range("H134").Value = 1867.5
$ws.Range("I134").Value = 1867.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5602.5
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -532.5

# CUL row 138 (Leve Item ID 44105)
$ws.Range("H138").Value = 1464.6364
$ws.Range("I138").Value = 1091.375
$ws.Range("J138").Value = 2460
$ws.Range("K138").Value = 3274.125
$ws.Range("L138").Value = 7380
$ws.Range("M138").Value = 1865.875
$ws.Range("N138").Value = -17660

$ws = $wb.Worksheets.Item("GSM")
# GSM row 6 (Leve Item ID 2108)
$ws.Range("H6").Value = 3000
$ws.Range("I6").Value = 3000
$ws.Range("K6").Value = 3000
$ws.Range("M6").Value = -2887

# GSM row 7 (Leve Item ID 4197)
$ws.Range("H7").Value = 5004875.5
$ws.Range("J7").Value = 2868287
$ws.Range("L7").Value = 2868287
$ws.Range("N7").Value = -2868511

# GSM row 8 (Leve Item ID 4197)
$ws.Range("H8").Value = 5004875.5
$ws.Range("J8").Value = 2868287
$ws.Range("L8").Value = 2868287
$ws.Range("N8").Value = -2868565

# GSM row 16 (Leve Item ID 2108)
$ws.Range("H16").Value = 3000
$ws.Range("I16").Value = 3000
$ws.Range("K16").Value = 3000
$ws.Range("M16").Value = -2750

# GSM row 104 (Leve Item ID 18666)
$ws.Range("H104").Value = 40000
$ws.Range("J104").Value = 40000
$ws.Range("L104").Value = 40000
$ws.Range("N104").Value = -46988

# GSM row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 4183.483
$ws.Range("I132").Value = 4977
$ws.Range("J132").Value = 3623.353
$ws.Range("K132").Value = 14931
$ws.Range("L132").Value = 10870.059
$ws.Range("M132").Value = -12401
$ws.Range("N132").Value = -15930.059

$ws = $wb.Worksheets.Item("LTW")
# LTW row 26 (Leve Item ID 3559)
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").ClearContents()
$ws.Range("N26").Value = 0

# LTW row 45 (Leve Item ID 3851)
$ws.Range("H45").Value = 7510
$ws.Range("I45").Value = 6020
$ws.Range("K45").Value = 6020
$ws.Range("M45").Value = -5613

# LTW row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 1373.2322
$ws.Range("I46").Value = 970.8542
$ws.Range("J46").Value = 3787.5
$ws.Range("K46").Value = 970.8542
$ws.Range("L46").Value = 3787.5
$ws.Range("M46").Value = -782.8542
$ws.Range("N46").Value = -4163.5

# LTW row 124 (Leve Item ID 34264)
$ws.Range("H124").Value = 29000
$ws.Range("J124").Value = 29000
$ws.Range("L124").Value = 29000
$ws.Range("N124").Value = -38820

# LTW row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 3973.1765
$ws.Range("I132").Value = 2693
$ws.Range("K132").Value = 8079
$ws.Range("M132").Value = -5549

$ws = $wb.Worksheets.Item("WVR")
# WVR row 56 (Leve Item ID 10912)
$ws.Range("H56").Value = 14771.333
$ws.Range("J56").Value = 19657
$ws.Range("L56").Value = 19657
$ws.Range("N56").Value = -21085
